$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full post-edit values for columns B:K, rows 2-24 ("ifoCAST full series evaluation").
# The underlying data block shifts one column to the left (an extra forecast-horizon
# column of new data is appended for rows 2-14; rows 15-24 simply lose their final,
# previously truncated column). $null clears/removes a cell entirely.
$allRows = @(
    @(2.489322709745834, 10.35911539074583, -7.969890609254166, 0.1701913907458341, 1.958500390745834, -1.219618609254166, -1.514947609254166, 0.9994973907458341, 0.06977539074583408, 0.3047833907458341),
    @(10.04419022232486, -8.284815777675137, -0.1447337776751375, 1.643575222324863, -1.534543777675137, -1.829872777675137, 0.6845722223248626, -0.2451497776751375, -0.01014177767513746, -0.3172777776751375),
    @(-9.341190737857668, -1.201108737857667, 0.5872002621423329, -2.590918737857667, -2.886247737857667, -0.371802737857667, -1.301524737857667, -1.066516737857667, -1.373652737857667, -0.862976737857667),
    @(-0.3304623618867365, 1.457846638113264, -1.720272361886736, -2.015601361886736, 0.4988436381132635, -0.4308783618867365, -0.1958703618867365, -0.5030063618867365, 0.007669638113263488, -0.2796973618867365),
    @(1.44719542178338, -1.73092357821662, -2.02625257821662, 0.48819242178338, -0.44152957821662, -0.20652157821662, -0.51365757821662, -0.002981578216620012, -0.29034857821662, -0.26234957821662),
    @(-1.839538110432761, -2.134867110432761, 0.3795778895672392, -0.5501441104327608, -0.3151361104327608, -0.6222721104327609, -0.1115961104327608, -0.3989631104327608, -0.3709641104327608, -0.2783161104327608),
    @(-1.926765389465022, 0.5876796105349782, -0.3420423894650218, -0.1070343894650218, -0.4141703894650218, 0.09650561053497819, -0.1908613894650218, -0.1628623894650218, -0.0702143894650218, 0.1427616105349782),
    @(0.7041777809087787, -0.2255442190912214, 0.009463780908778634, -0.2976722190912214, 0.2130037809087786, -0.07436321909122137, -0.04636421909122137, 0.04628378090877863, 0.2592597809087786, -0.2414002190912214),
    @(-0.1800508050940863, 0.05495719490591373, -0.2521788050940862, 0.2584971949059137, -0.02886980509408627, -0.0008708050940862688, 0.09177719490591374, 0.3047531949059137, -0.1959068050940863, 0.09385019490591373),
    @(0.03740307354594928, -0.2697329264540507, 0.2409430735459493, -0.04642392645405072, -0.01842492645405072, 0.07422307354594927, 0.2871990735459493, -0.2134609264540507, 0.07629607354594928, -0.1815959264540507),
    @(-0.1816256075030445, 0.3290503924969555, 0.0416833924969555, 0.06968239249695551, 0.1623303924969555, 0.3753063924969555, -0.1253536075030445, 0.1644033924969555, -0.0934886075030445, 0.2888923924969555),
    @(0.5076250259847391, 0.2202580259847391, 0.2482570259847391, 0.3409050259847391, 0.5538810259847391, 0.0532210259847391, 0.3429780259847391, 0.0850860259847391, 0.4674670259847391, 0.05618993198473912),
    @(1.132285253737389, 1.160284253737389, 1.252932253737389, 1.465908253737389, 0.9652482537373891, 1.255005253737389, 0.9971132537373891, 1.379494253737389, 0.9682171597373891, 1.255173253737389),
    @(0.2163100177716323, 0.3089580177716323, 0.5219340177716323, 0.02127401777163229, 0.3110310177716323, 0.0531390177716323, 0.4355200177716323, 0.02424292377163231, 0.3111990177716323, $null),
    @(0.3129259152847414, 0.5259019152847414, 0.0252419152847414, 0.3149989152847414, 0.05710691528474141, 0.4394879152847414, 0.02821082128474142, 0.3151669152847414, $null, $null),
    @(0.6615420054549828, 0.1608820054549828, 0.4506390054549828, 0.1927470054549828, 0.5751280054549828, 0.1638509114549828, 0.4508070054549828, $null, $null, $null),
    @(-0.07992417640068192, 0.2098328235993181, -0.04805917640068191, 0.3343218235993181, -0.0769552704006819, 0.2100008235993181, $null, $null, $null, $null),
    @(0.1645016083472217, -0.09339039165277828, 0.2889906083472217, -0.1222864856527783, 0.1646696083472217, $null, $null, $null, $null, $null),
    @(-0.1602518135888321, 0.2221291864111679, -0.1891479075888321, 0.09780818641116784, $null, $null, $null, $null, $null, $null),
    @(0.1925429840544512, -0.2187341099455487, 0.06822198405445125, $null, $null, $null, $null, $null, $null, $null),
    @(-0.2368847666751896, 0.05007132732481036, $null, $null, $null, $null, $null, $null, $null, $null),
    @(0.03203950307611103, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

for ($r = 2; $r -le 24; $r++) {
    $vals = $allRows[$r - 2]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]   # column B = index 2
    }
}
